$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "62.372.05"
Set-TextCell 2 5 "  -1.72%  "
Set-TextCell 3 4 "3.012.28"
Set-TextCell 3 5 "  -2.19%  "
Set-TextCell 4 5 "  -0.13%  "
Set-TextCell 5 4 "577.94"
Set-TextCell 5 5 "  -1.90%  "
Set-TextCell 6 4 "148.79"
Set-TextCell 6 5 "  -3.08%  "
Set-TextCell 7 5 "  -0.06%  "
Set-TextCell 8 4 "0.523"
Set-TextCell 8 5 "  -3.34%  "
Set-TextCell 9 4 "3.009.47"
Set-TextCell 9 5 "  -1.93%  "
Set-TextCell 10 5 "  -4.32%  "
Set-TextCell 11 4 "5.68"
Set-TextCell 11 5 "  -1.92%  "
Set-TextCell 12 4 "0.442"
Set-TextCell 12 5 "  -2.69%  "
Set-TextCell 13 5 "  -3.97%  "
Set-TextCell 14 4 "35.36"
Set-TextCell 14 5 "  -5.03%  "
Set-TextCell 15 4 "0.120"
Set-TextCell 15 5 "  +1.49%  "
Set-TextCell 16 4 "3.513.78"
Set-TextCell 16 5 "  -2.13%  "
Set-TextCell 17 2 "WrappedBTC"
Set-TextCell 17 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 17 4 "62.409.53"
Set-TextCell 17 5 "  -1.62%  "
Set-TextCell 18 2 "Polkadot"
Set-TextCell 18 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 18 4 "7.00"
Set-TextCell 18 5 "  -2.05%  "
Set-TextCell 19 4 "3.016.13"
Set-TextCell 19 5 "  -1.99%  "
Set-TextCell 20 4 "470.13"
Set-TextCell 20 5 "  -0.81%  "
Set-TextCell 21 4 "14.01"
Set-TextCell 21 5 "  -3.61%  "
Set-TextCell 22 4 "0.691"
Set-TextCell 22 5 "  -2.87%  "
Set-TextCell 23 4 "7.39"
Set-TextCell 23 5 "  -1.06%  "
Set-TextCell 24 4 "2.34"
Set-TextCell 24 5 "  -1.77%  "
Set-TextCell 25 4 "80.71"
Set-TextCell 25 5 "  -0.31%  "
Set-TextCell 26 4 "12.41"
Set-TextCell 26 5 "  -3.29%  "
Set-TextCell 27 4 "10.43"
Set-TextCell 27 5 "  +4.82%  "
Set-TextCell 28 4 "0.999"
Set-TextCell 28 5 "  +0.08%  "
Set-TextCell 29 4 "1.00"
Set-TextCell 29 5 "  -0.02%  "
Set-TextCell 30 4 "7.17"
Set-TextCell 30 5 "  -1.72%  "
Set-TextCell 31 4 "2.61"
Set-TextCell 31 5 "  -2.32%  "
Set-TextCell 32 4 "2.16"
Set-TextCell 32 5 "  -0.49%  "
Set-TextCell 33 4 "27.10"
Set-TextCell 33 5 "  +0.03%  "
Set-TextCell 34 4 "0.108"
Set-TextCell 34 5 "  -4.50%  "
Set-TextCell 35 5 "  -1.16%  "
Set-TextCell 36 4 "0.0₃0793"
Set-TextCell 36 5 "  -5.89%  "
Set-TextCell 37 4 "5.79"
Set-TextCell 37 5 "  -3.99%  "
Set-TextCell 38 4 "2.15"
Set-TextCell 38 5 "  -2.21%  "
Set-TextCell 39 5 "  -9.94%  "
Set-TextCell 40 4 "50.07"
Set-TextCell 40 5 "  -0.97%  "
Set-TextCell 41 4 "8.98"
Set-TextCell 41 5 "  -2.56%  "
Set-TextCell 42 4 "419.08"
Set-TextCell 42 5 "  -4.86%  "
Set-TextCell 43 5 "  +2.10%  "
Set-TextCell 44 4 "0.279"
Set-TextCell 44 5 "  -0.98%  "
Set-TextCell 45 4 "2.801.51"
Set-TextCell 45 5 "  +0.33%  "
Set-TextCell 46 4 "0.0355"
Set-TextCell 46 5 "  -1.40%  "
Set-TextCell 47 4 "38.05"
Set-TextCell 47 5 "  -4.26%  "
Set-TextCell 48 4 "127.33"
Set-TextCell 48 5 "  -2.62%  "
Set-TextCell 49 5 "  +0.01%  "
Set-TextCell 50 4 "24.73"
Set-TextCell 50 5 "  -3.21%  "
Set-TextCell 51 5 "  -1.46%  "
